$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" (copy the header style from an existing
#     header cell, e.g. H1, so the bold/centered/bordered style carries over) ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-49 for columns I (I0) and J (IF) ---
$data = @(
  @(6,7),
  @(7,7),
  @(6,8),
  @(11,11),
  @(6,6),
  @(8,8),
  @(4,6),
  @(1,3),
  @(4,6),
  @(7,8),
  @(7,8),
  @(5,5),
  @(6,8),
  @(6,7),
  @(8,8),
  @(7,7),
  @(6,6),
  @(7,9),
  @(9,9),
  @(8,8),
  @(6,7),
  @(5,6),
  @(7,8),
  @(6,7),
  @(5,7),
  @(7,7),
  @(6,9),
  @(7,9),
  @(5,7),
  @(6,8),
  @(5,6),
  @(8,9),
  @(5,7),
  @(7,9),
  @(7,9),
  @(5,8),
  @(3,5),
  @(7,8),
  @(7,9),
  @(7,9),
  @(4,5),
  @(5,6),
  @(7,8),
  @(4,6),
  @(6,8),
  @(5,7),
  @(4,5),
  @(5,5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
